$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2:G13").Value = "2000 : 1"
$ws.Range("H2:H13").Value = "500 : 1"
$ws.Range("G2:H13").WrapText = $true

$ws.Range("G2").Select()
